$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "66.276.15"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +6.89%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.008.26"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +3.78%  "
$ws.Range("E4").Value = "  -0.10%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "582.57"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +3.01%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "162.99"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +13.59%  "
$ws.Range("E7").Value = "  -0.07%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "3.003.73"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +3.84%  "
$ws.Range("E9").Value = "  +3.68%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "6.63"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -4.36%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.156"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +5.13%  "
$ws.Range("E12").Value = "  +5.59%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.0000256"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +7.63%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "34.63"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +6.82%  "
$ws.Range("E15").Value = "  -0.83%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "66.267.12"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +6.89%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "3.507.14"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +3.76%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "6.93"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +5.67%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "3.010.72"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +3.46%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "454.72"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +6.04%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "13.84"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +6.38%  "
$ws.Range("E22").Value = "  +4.74%  "
$ws.Range("E23").Value = "  +7.67%  "
$ws.Range("E24").Value = "  +4.83%  "
$ws.Range("E25").Value = "  +15.80%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "12.27"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +3.07%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "10.48"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +6.24%  "
$ws.Range("E28").Value = "  +0.03%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "8.17"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +18.07%  "
$ws.Range("E30").Value = "  +20.63%  "
$ws.Range("E31").Value = "  -4.49%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "2.60"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +5.01%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "27.28"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +6.67%  "
$ws.Range("E34").Value = "  +5.36%  "
$ws.Range("E35").Value = "  -0.10%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.993"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +3.99%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "5.81"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +8.39%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "2.17"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +15.48%  "
$ws.Range("E39").Value = "  +3.96%  "
$ws.Range("E40").Value = "  +2.34%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.309"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +17.28%  "
$ws.Range("E42").Value = "  +8.50%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "43.85"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +7.99%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "8.42"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +4.20%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "396.78"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +14.58%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.0360"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +7.74%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "2.793.52"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +3.10%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "134.26"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +1.13%  "
$ws.Range("E49").Value = "  +0.00%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "23.90"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +12.29%  "
$ws.Range("E51").Value = "  +4.88%  "
